$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "313.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.15%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "8"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.43%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "8"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.135"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.24%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "8"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08107"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.11%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "8"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.477"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.50%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "8"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.961"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.49%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "8"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.299"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.75%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "8"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9380"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.03%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "8"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1330"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.59%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "8"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1963"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.63%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "8"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09072"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.56%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "8"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03483"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.81%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "8"

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.11%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "8"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001404"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.02%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "8"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005973"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.07%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "8"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.555"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-8.89%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "8"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.193"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.89%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "8"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3458"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.09%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "8"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1290"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.72%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "8"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.007"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.49%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "8"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2490"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.60%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "8"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04370"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.25%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "8"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001246"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.27%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "8"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004734"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.38%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "8"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003892"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "199.17%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "8"

# Row 27
$ws.Range("B27").Value = "Spectre.aiUtilityToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "--"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "--%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "8"

# Row 28
$ws.Range("B28").Value = "LegolasExchange"
$ws.Range("C28").Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "8"

# Row 29
$ws.Range("B29").Value = "BitZToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "8"

# Row 30
$ws.Range("B30").Value = "Birake"
$ws.Range("C30").Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "8"

# Row 31
$ws.Range("B31").Value = "NashExchange"
$ws.Range("C31").Value = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "8"

# Row 32
$ws.Range("B32").Value = "AAXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "8"

# Row 33
$ws.Range("B33").Value = "CenX"
$ws.Range("C33").Value = "https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "8"

# Row 34
$ws.Range("B34").Value = "BNIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "8"

# Row 35
$ws.Range("B35").Value = "UpBots"
$ws.Range("C35").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "8"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "8"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "8"

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "8"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02213"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.93%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "8"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05215"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.16%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "8"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007619"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.24%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "8"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01033"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.03%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "8"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1395"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.38%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "8"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.43%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "8"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009097"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.73%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "8"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006604"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.15%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "8"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "8"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003013"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "16.53%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "8"

# Row 49
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "8"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "8"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "8"
